$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Preserve the full recipe list: duplicate the "recipes" sheet to the end
#    of the workbook and rename the copy "all_recipes" (so the master list of
#    every scraped recipe URL survives).
# ---------------------------------------------------------------------------
$recipes = $wb.Worksheets.Item("recipes")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$recipes.Copy($null, $lastSheet)
$allRecipes = $wb.Worksheets.Item($wb.Worksheets.Count)
$allRecipes.Name = "all_recipes"

# ---------------------------------------------------------------------------
# 2. Stop auto-picking a random set of recipes: the working "recipes" sheet
#    now keeps just two hand-picked entries (what used to be row 10 moves up
#    to row 2, what used to be row 2 moves to row 3) and the remaining rows
#    (4-10) are blanked out along with their hyperlinks.
# ---------------------------------------------------------------------------
$oldA2 = $recipes.Range("A2").Value2
$oldA10 = $recipes.Range("A10").Value2

$hlA2Target = $null
$hlA10Target = $null
foreach ($hl in @($recipes.Hyperlinks)) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') { $hlA2Target = $hl.Address }
    if ($addr -eq '$A$10') { $hlA10Target = $hl.Address }
}

foreach ($hl in @($recipes.Hyperlinks)) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.Address = $hlA10Target
    }
}
$recipes.Range("A2").Value = $oldA10

foreach ($hl in @($recipes.Hyperlinks)) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$3') {
        $hl.Address = $hlA2Target
    }
}
$recipes.Range("A3").Value = $oldA2

# Clear rows 4-10 and drop their hyperlinks entirely.
$recipes.Range("A4:A10").ClearContents()
$toDrop = @('$A$4', '$A$5', '$A$6', '$A$7', '$A$8', '$A$9', '$A$10')
$keepGoing = $true
while ($keepGoing) {
    $keepGoing = $false
    foreach ($hl in @($recipes.Hyperlinks)) {
        $addr = $hl.Range.Address()
        if ($toDrop -contains $addr) {
            $hl.Delete()
            $keepGoing = $true
            break
        }
    }
}

$recipes.Range("A4").Select() | Out-Null

# ---------------------------------------------------------------------------
# 3. "anything premade" is now a recognized additional item.
# ---------------------------------------------------------------------------
$allItems = $wb.Worksheets.Item("all_additional_items")
$allItems.Range("A22").Value = "anything premade"
